$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 3.5
$ws.Range("G2").Value = 4.5
$ws.Range("H2").Value = 2.18
$ws.Range("I2").Value = 2.56
$ws.Range("K2").Value = 3.8
$ws.Range("L2").Value = 1.46
$ws.Range("M2").Value = 1.11
$ws.Range("N2").Value = 2.54
$ws.Range("O2").Value = 1.49
$ws.Range("P2").Value = 1.52
$ws.Range("Q2").Value = 2.44
$ws.Range("R2").Value = 1.19
$ws.Range("S2").Value = 4.5
$ws.Range("T2").Value = 1.89
$ws.Range("U2").Value = 1.66
$ws.Range("V2").Value = 1.64
$ws.Range("W2").Value = 1.3

# Row 3
$ws.Range("H3").Value = 2.7
$ws.Range("I3").Value = 3.85
$ws.Range("K3").Value = 5.1
$ws.Range("L3").Value = 1.52
$ws.Range("V3").Value = 1.35

# Row 4
$ws.Range("H4").Value = 2.76
$ws.Range("I4").Value = 2.94
$ws.Range("K4").Value = 2.94
$ws.Range("M4").Value = 1.2
$ws.Range("N4").Value = 2
$ws.Range("O4").Value = 1.84
$ws.Range("P4").Value = 1.33
$ws.Range("T4").Value = 2.66
$ws.Range("U4").Value = 1.44
$ws.Range("W4").Value = 1.4
$ws.Range("X4").Value = 5.6
$ws.Range("Y4").Value = 7
$ws.Range("Z4").Value = 16.5
$ws.Range("AA4").Value = 65
$ws.Range("AC4").Value = 8
$ws.Range("AD4").Value = 17.5
$ws.Range("AE4").Value = 70
$ws.Range("AF4").Value = 21
$ws.Range("AG4").Value = 19
$ws.Range("AH4").Value = 38
$ws.Range("AJ4").Value = 90
$ws.Range("AK4").Value = 85
$ws.Range("AL4").Value = 170
$ws.Range("AM4").Value = 440
$ws.Range("AN4").Value = 140
$ws.Range("AO4").Value = 110

# Row 6
$ws.Range("T6").Value = 2.12

# Row 7
$ws.Range("J7").Value = 3
$ws.Range("Q7").Value = 2.76
$ws.Range("T7").Value = 2.02
